$d = $word.ActiveDocument

# 1. "februar 2020 kl." -> "februar 2020, kl." (insert comma after the year)
$d.Content.Find.Execute("2020 kl.", $true, $false, $false, $false, $false, $true, 1, $false, "2020, kl.", 2)

# 2. Agenda list item: "Valg af kasserer frem til generalforsamlingen i 2020" -> "Valg af revisor frem til generalforsamlingen i 2020"
$d.Content.Find.Execute("Valg af kasserer frem", $true, $false, $false, $false, $false, $true, 1, $false, "Valg af revisor frem", 2)

# 3. "Jesper Kastoft Bertelsen og Tobias Evars Lauridsen blev valgt som stemmetællere." -> fix "Evars" -> "Evar"
$d.Content.Find.Execute("Tobias Evars Lauridsen blev valgt som stemmetællere.", $true, $false, $false, $false, $false, $true, 1, $false, "Tobias Evar Lauridsen blev valgt som stemmetællere.", 2)

# 4. "Tobias Evars Lauridsen med 4 stemmer" -> fix "Evars" -> "Evar"
$d.Content.Find.Execute("Tobias Evars Lauridsen med 4 stemmer", $true, $false, $false, $false, $false, $true, 1, $false, "Tobias Evar Lauridsen med 4 stemmer", 2)

# 5. "Ad 9. Valg af kasserer" heading -> "Ad 9. Valg af revisor"
$d.Content.Find.Execute("Ad 9. Valg af kasserer", $true, $false, $false, $false, $false, $true, 1, $false, "Ad 9. Valg af revisor", 2)

# 6. Body of Ad 9 paragraph: replace the whole sentence about who was elected
$d.Content.Find.Execute("Tobias Evars Lauridsen blev valgt og påtager sig rollen som kasserer.", $true, $false, $false, $false, $false, $true, 1, $false, "Morten Marthendal Lond påtager sig rollen revisor for BSides Aarhus.", 2)

# 7. Move the "_GoBack" bookmark from the "forslag" paragraph to just before
#    "revisor for BSides Aarhus." at the end of the document (matching the
#    author's final cursor position after editing the Ad 9 paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$target = $d.Content
$target.Find.Execute("rollen ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)
$d.Bookmarks.Add("_GoBack", $target)
